# lecture03.pptx edit — re-saved deck:
#   * title slide subtitle updated ("January 14, 2022" -> "Fall 2022" /
#     "The university of mount union")
#   * every cached "datetimeFigureOut" Date placeholder (slide master +
#     all slide layouts) refreshed to the new save date, 10/17/2022

$p = $ppt.ActivePresentation

# --- 1. Title slide (slide 1) subtitle text -------------------------------
$slide1 = $p.Slides.Item(1)
for ($i = 1; $i -le $slide1.Shapes.Count; $i++) {
    $shp = $slide1.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.Name -like "Subtitle*") {
        # `r` starts a new <a:p> paragraph, matching the two-line subtitle.
        $shp.TextFrame.TextRange.Text = "Fall 2022`rThe university of mount union"
    }
}

# --- 2. Refresh the Date placeholder on the master + every layout --------
function Update-DatePlaceholder($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)
        if ($shp.HasTextFrame -and $shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = "10/17/2022"
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

for ($k = 1; $k -le $master.CustomLayouts.Count; $k++) {
    $layout = $master.CustomLayouts.Item($k)
    Update-DatePlaceholder $layout.Shapes
}
